$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column S (year 2022) data, mirroring the styles used in column R
$ws.Cells.Item(4, 19).Value = 2022
$ws.Cells.Item(5, 19).Value = 4.9538761752705343
$ws.Cells.Item(6, 19).Value = 11.304954640614097
$ws.Cells.Item(7, 19).Value = 5.1593323216995444
$ws.Cells.Item(8, 19).Value = 13.687943262411348
$ws.Cells.Item(9, 19).Value = 10.22864019253911
$ws.Cells.Item(10, 19).Value = 9.1213700670141478
$ws.Cells.Item(11, 19).Value = 3.1335149863760217
$ws.Cells.Item(12, 19).Value = 2.872905173311127
$ws.Cells.Item(13, 19).Value = 3.527842284697861
$ws.Cells.Item(14, 19).Value = 5.0305321314335565

# Copy formatting from column R to column S for rows 4-14
$ws.Range("R4:R14").Copy()
$ws.Range("S4:S14").PasteSpecial(-4122)

# Update the selection to match the recorded state after the edit
$ws.Range("T6").Select()
